$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H69").Value = 15777.944
$ws_ALC.Range("I69").Value = 8337.666999999999
$ws_ALC.Range("J69").Value = 17266
$ws_ALC.Range("K69").Value = 25013.001
$ws_ALC.Range("L69").Value = 51798
$ws_ALC.Range("M69").Value = -24139.001
$ws_ALC.Range("N69").Value = -53546

$ws_ALC.Range("H70").Value = 1431

$ws_ALC.Range("H72").Value = 15777.944
$ws_ALC.Range("I72").Value = 8337.666999999999
$ws_ALC.Range("J72").Value = 17266
$ws_ALC.Range("K72").Value = 75039.003
$ws_ALC.Range("L72").Value = 155394
$ws_ALC.Range("M72").Value = -70671.003
$ws_ALC.Range("N72").Value = -164130

$ws_ALC.Range("H73").Value = 1431

$ws_ALC.Range("H80").Value = 1730.7916
$ws_ALC.Range("I80").Value = 1589.125
$ws_ALC.Range("J80").Value = 1801.625
$ws_ALC.Range("K80").Value = 4767.375
$ws_ALC.Range("L80").Value = 5404.875
$ws_ALC.Range("M80").Value = -3769.375
$ws_ALC.Range("N80").Value = -7400.875

$ws_ALC.Range("H82").Value = 904.1539
$ws_ALC.Range("I82").Value = 904.1539
$ws_ALC.Range("K82").Value = 2712.4617
$ws_ALC.Range("M82").Value = -2306.4617

$ws_ALC.Range("H83").Value = 1730.7916
$ws_ALC.Range("I83").Value = 1589.125
$ws_ALC.Range("J83").Value = 1801.625
$ws_ALC.Range("K83").Value = 14302.125
$ws_ALC.Range("L83").Value = 16214.625
$ws_ALC.Range("M83").Value = -9310.125
$ws_ALC.Range("N83").Value = -26198.625

$ws_ALC.Range("H85").Value = 904.1539
$ws_ALC.Range("I85").Value = 904.1539
$ws_ALC.Range("K85").Value = 2712.4617
$ws_ALC.Range("M85").Value = -1308.4617

$ws_ALC.Range("H88").Value = 6044.5884
$ws_ALC.Range("I88").Value = 4956
$ws_ALC.Range("J88").Value = 6498.1665
$ws_ALC.Range("K88").Value = 4956
$ws_ALC.Range("L88").Value = 6498.1665
$ws_ALC.Range("M88").Value = -4550
$ws_ALC.Range("N88").Value = -7310.1665

$ws_ALC.Range("H91").Value = 6044.5884
$ws_ALC.Range("I91").Value = 4956
$ws_ALC.Range("J91").Value = 6498.1665
$ws_ALC.Range("K91").Value = 4956
$ws_ALC.Range("L91").Value = 6498.1665
$ws_ALC.Range("M91").Value = -3552
$ws_ALC.Range("N91").Value = -9306.166499999999

$ws_ALC.Range("H131").Value = 250003890
$ws_ALC.Range("I131").Value = 333337700
$ws_ALC.Range("K131").Value = 1000013100
$ws_ALC.Range("M131").Value = -1000008060

$ws_ARM.Range("H32").Value = 4652.908
$ws_ARM.Range("I32").Value = 3048.782
$ws_ARM.Range("J32").Value = 18555.334
$ws_ARM.Range("K32").Value = 3048.782
$ws_ARM.Range("L32").Value = 18555.334
$ws_ARM.Range("M32").Value = -2761.782
$ws_ARM.Range("N32").Value = -19129.334

$ws_ARM.Range("H45").Value = 1872.5
$ws_ARM.Range("I45").Value = 1376.8889
$ws_ARM.Range("K45").Value = 1376.8889
$ws_ARM.Range("M45").Value = -999.8888999999999

$ws_ARM.Range("H88").Value = 2388.3635
$ws_ARM.Range("I88").Value = 2162
$ws_ARM.Range("J88").Value = 2660
$ws_ARM.Range("K88").Value = 2162
$ws_ARM.Range("L88").Value = 2660
$ws_ARM.Range("M88").Value = -1756
$ws_ARM.Range("N88").Value = -3472

$ws_ARM.Range("H91").Value = 2388.3635
$ws_ARM.Range("I91").Value = 2162
$ws_ARM.Range("J91").Value = 2660
$ws_ARM.Range("K91").Value = 2162
$ws_ARM.Range("L91").Value = 2660
$ws_ARM.Range("M91").Value = -758
$ws_ARM.Range("N91").Value = -5468

$ws_BSM.Range("H86").Value = 1457.4147
$ws_BSM.Range("I86").Value = 1368.9584
$ws_BSM.Range("K86").Value = 1368.9584
$ws_BSM.Range("M86").Value = -245.9584

$ws_BSM.Range("H89").Value = 1457.4147
$ws_BSM.Range("I89").Value = 1368.9584
$ws_BSM.Range("K89").Value = 6844.791999999999
$ws_BSM.Range("M89").Value = -1228.791999999999

$ws_BSM.Range("H107").Value = 2344.353
$ws_BSM.Range("I107").Value = 2236.75
$ws_BSM.Range("J107").Value = 2602.6
$ws_BSM.Range("K107").Value = 2236.75
$ws_BSM.Range("L107").Value = 2602.6
$ws_BSM.Range("M107").Value = -316.75
$ws_BSM.Range("N107").Value = -6442.6

$ws_CRP.Range("H31").Value = 1377.5
$ws_CRP.Range("I31").Value = 842.6429000000001
$ws_CRP.Range("J31").Value = 2781.5
$ws_CRP.Range("K31").Value = 842.6429000000001
$ws_CRP.Range("L31").Value = 2781.5
$ws_CRP.Range("M31").Value = -547.6429000000001
$ws_CRP.Range("N31").Value = -3371.5

$ws_CRP.Range("H34").Value = 1377.5
$ws_CRP.Range("I34").Value = 842.6429000000001
$ws_CRP.Range("J34").Value = 2781.5
$ws_CRP.Range("K34").Value = 842.6429000000001
$ws_CRP.Range("L34").Value = 2781.5
$ws_CRP.Range("M34").Value = -640.6429000000001
$ws_CRP.Range("N34").Value = -3185.5

$ws_CRP.Range("H58").Value = 4328.2188
$ws_CRP.Range("I58").Value = 5163.64
$ws_CRP.Range("J58").Value = 1344.5714
$ws_CRP.Range("K58").Value = 5163.64
$ws_CRP.Range("L58").Value = 1344.5714
$ws_CRP.Range("M58").Value = -4960.64
$ws_CRP.Range("N58").Value = -1750.5714

$ws_CRP.Range("H99").Value = 3808.6943
$ws_CRP.Range("I99").Value = 3602.2307
$ws_CRP.Range("J99").Value = 4345.5
$ws_CRP.Range("K99").Value = 3602.2307
$ws_CRP.Range("L99").Value = 4345.5
$ws_CRP.Range("M99").Value = -2104.2307
$ws_CRP.Range("N99").Value = -7341.5

$ws_CRP.Range("H126").Value = 3808.6943
$ws_CRP.Range("I126").Value = 3602.2307
$ws_CRP.Range("J126").Value = 4345.5
$ws_CRP.Range("K126").Value = 10806.6921
$ws_CRP.Range("L126").Value = 13036.5
$ws_CRP.Range("M126").Value = -8336.6921
$ws_CRP.Range("N126").Value = -17976.5

$ws_CRP.Range("H132").Value = 1002940.4
$ws_CRP.Range("I132").Value = 2046.3103
$ws_CRP.Range("K132").Value = 6138.9309
$ws_CRP.Range("M132").Value = -3608.9309

$ws_CRP.Range("H136").Value = 4328.2188
$ws_CRP.Range("I136").Value = 5163.64
$ws_CRP.Range("J136").Value = 1344.5714
$ws_CRP.Range("K136").Value = 15490.92
$ws_CRP.Range("L136").Value = 4033.7142
$ws_CRP.Range("M136").Value = -12940.92
$ws_CRP.Range("N136").Value = -9133.7142

$ws_CUL.Range("H131").Value = 913.45
$ws_CUL.Range("I131").Value = 500
$ws_CUL.Range("J131").Value = 917.6263
$ws_CUL.Range("K131").Value = 1500
$ws_CUL.Range("L131").Value = 2752.8789
$ws_CUL.Range("M131").Value = 3540
$ws_CUL.Range("N131").Value = -12832.8789

$ws_GSM.Range("H80").Value = 2424.95
$ws_GSM.Range("I80").Value = 2323.4707
$ws_GSM.Range("J80").Value = 3000
$ws_GSM.Range("K80").Value = 2323.4707
$ws_GSM.Range("L80").Value = 3000
$ws_GSM.Range("M80").Value = -1325.4707
$ws_GSM.Range("N80").Value = -4996

$ws_GSM.Range("H83").Value = 2424.95
$ws_GSM.Range("I83").Value = 2323.4707
$ws_GSM.Range("J83").Value = 3000
$ws_GSM.Range("K83").Value = 11617.3535
$ws_GSM.Range("L83").Value = 15000
$ws_GSM.Range("M83").Value = -6625.353499999999
$ws_GSM.Range("N83").Value = -24984

$ws_LTW.Range("H68").Value = 8535.235000000001
$ws_LTW.Range("I68").Value = 14150.5
$ws_LTW.Range("J68").Value = 3543.889
$ws_LTW.Range("K68").Value = 14150.5
$ws_LTW.Range("L68").Value = 3543.889
$ws_LTW.Range("M68").Value = -13401.5
$ws_LTW.Range("N68").Value = -5041.889

$ws_LTW.Range("H71").Value = 8535.235000000001
$ws_LTW.Range("I71").Value = 14150.5
$ws_LTW.Range("J71").Value = 3543.889
$ws_LTW.Range("K71").Value = 70752.5
$ws_LTW.Range("L71").Value = 17719.445
$ws_LTW.Range("M71").Value = -67008.5
$ws_LTW.Range("N71").Value = -25207.445

$ws_LTW.Range("H82").Value = 1267.2593
$ws_LTW.Range("I82").Value = 1447.6666
$ws_LTW.Range("J82").Value = 1122.9333
$ws_LTW.Range("K82").Value = 1447.6666
$ws_LTW.Range("L82").Value = 1122.9333
$ws_LTW.Range("M82").Value = -1086.6666
$ws_LTW.Range("N82").Value = -1844.9333

$ws_LTW.Range("H85").Value = 1267.2593
$ws_LTW.Range("I85").Value = 1447.6666
$ws_LTW.Range("J85").Value = 1122.9333
$ws_LTW.Range("K85").Value = 1447.6666
$ws_LTW.Range("L85").Value = 1122.9333
$ws_LTW.Range("M85").Value = -199.6666
$ws_LTW.Range("N85").Value = -3618.9333

$ws_WVR.Range("H81").Value = 1160.4
$ws_WVR.Range("I81").Value = 1160.4
$ws_WVR.Range("K81").Value = 2320.8
$ws_WVR.Range("M81").Value = -1259.8

$ws_WVR.Range("H84").Value = 1160.4
$ws_WVR.Range("I84").Value = 1160.4
$ws_WVR.Range("K84").Value = 11604
$ws_WVR.Range("M84").Value = -6300

$ws_WVR.Range("H132").Value = 2711.195
$ws_WVR.Range("I132").Value = 3084.7932
$ws_WVR.Range("K132").Value = 9254.3796
$ws_WVR.Range("M132").Value = -6724.3796
